$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block in rows 2-4 (3 columns wide)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 3

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 3

# Single-column trailer rows
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# Row 8 is intentionally left blank; old A3 value (8) now lives at A9
$ws.Range("A9").Value = 8
